# Update the cryptocurrency price/volume/coin data on sheet1 to reflect
# the latest snapshot (coin rows 2-51), per commit:
# 'Updated cryptos list on Fri May 26 22:10:23 UTC 2023 with GitHub Actions'

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '26.870.84'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '1.841.53'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.52'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4705'
$ws.Range('E7').Value = '  +3.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3661'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07152'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  +3.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.56'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07670'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').Value = '1.861.19'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.287'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.396'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.27'
$ws.Range('E16').Value = '  +2.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008632'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '26.909.06'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.46'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.008'
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.926'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.86'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.23'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.012'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.10'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.883'
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08820'
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.212'
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.178'
$ws.Range('E32').Value = '  +6.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7461'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.767'
$ws.Range('E34').Value = '  +1.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.480'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01941'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05207'
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.957'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5198'
$ws.Range('E40').Value = '  +1.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.975'
$ws.Range('E41').Value = '  +2.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1509'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.156'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.49'
$ws.Range('E44').Value = '  +5.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4698'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.007'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.73'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.598'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '66.00'
$ws.Range('E49').Value = '  +3.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06039'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8886'
$ws.Range('E51').Value = '  +5.10%  '
